# cambio de de debut sexual y nombres de ejes en español
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo "serorreversion" -> "seroreversion" in header row 8 (F8, G8)
$ws.Range("F8").Value = "Tasa de seroreversion"
$ws.Range("G8").Value = "Tasa de seroreversión Rhat"

# Update the selection to cover the whole second table (A8:G12)
$ws.Range("A8:G12").Select()
